$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns C and D to fit the new data
# (ColumnWidth is internally snapped to 1/6-character increments by this
# runtime, so 14.0 / 15.0 are the inputs that land closest to the target
# stored widths of 14.853482 / 15.853482.)
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 15.0

# Add three new rows (16-18) below the existing table, reusing the same
# alternating row styles (20/21 for row16 & row18, 22/23 for row17) that
# rows 12-15 already use.

$ws.Range("B16").Value = "Steam"
$ws.Range("C16").Value = "RogueLikeStore"
$ws.Range("D16").Value = "Hades"
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 60

$ws.Range("B17").Value = "Steam"
$ws.Range("C17").Value = "RogueLikeStore"
$ws.Range("D17").Value = "Dark Souls 3"
$ws.Range("E17").Value = 12
$ws.Range("F17").Value = 80

$ws.Range("B18").Value = "Steam"
$ws.Range("C18").Value = "RogueLikeStore"
$ws.Range("D18").Value = "Cult of The Lamb"
$ws.Range("E18").Value = 120
$ws.Range("F18").Value = 400

# Copy styling from existing rows (row14 has style pattern 20/21, row15 has 22/23)
$ws.Range("B14:G14").Copy()
$ws.Range("B16:G16").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B15:G15").Copy()
$ws.Range("B17:G17").PasteSpecial(-4122)

$ws.Range("B14:G14").Copy()
$ws.Range("B18:G18").PasteSpecial(-4122)

$excel.CutCopyMode = 0
